$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a price-column value while forcing text storage, since
# values like "393.88" or "1.00" would otherwise be auto-parsed as numbers
# by Excel, losing the trailing/leading zero formatting used in this sheet.
function Set-PriceText {
    param($cellRef, $text)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2 - Bitcoin
Set-PriceText "D2" "51.370.96"
$ws.Range("E2").Value = "  -0.51%  "

# Row 3 - Ethereum
Set-PriceText "D3" "3.074.86"
$ws.Range("E3").Value = "  +1.06%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-PriceText "D5" "393.88"
$ws.Range("E5").Value = "  +2.52%  "

# Row 6 - Solana
Set-PriceText "D6" "102.41"
$ws.Range("E6").Value = "  -0.38%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.68%  "

# Row 9 - Cardano
Set-PriceText "D9" "0.588"
$ws.Range("E9").Value = "  +0.39%  "

# Row 10 - Avalanche
Set-PriceText "D10" "37.48"
$ws.Range("E10").Value = "  +1.53%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.66%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  -1.33%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-PriceText "D13" "3.555.12"
$ws.Range("E13").Value = "  +0.85%  "

# Row 14 - Chainlink
Set-PriceText "D14" "18.62"
$ws.Range("E14").Value = "  -0.54%  "

# Row 15 - Polkadot
Set-PriceText "D15" "7.69"
$ws.Range("E15").Value = "  -0.69%  "

# Row 16 - Polygon
Set-PriceText "D16" "1.03"
$ws.Range("E16").Value = "  +5.01%  "

# Row 17 - WrappedEther
Set-PriceText "D17" "3.052.11"
$ws.Range("E17").Value = "  -0.23%  "

# Row 18 - Uniswap
Set-PriceText "D18" "10.55"
$ws.Range("E18").Value = "  -0.15%  "

# Row 19 - WrappedBTC
Set-PriceText "D19" "51.372.44"
$ws.Range("E19").Value = "  -0.58%  "

# Row 20 - ImmutableX
$ws.Range("E20").Value = "  +1.97%  "

# Row 21 - InternetComputer(DFINITY)
Set-PriceText "D21" "12.36"
$ws.Range("E21").Value = "  -0.91%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  -0.45%  "

# Row 23 - Litecoin
Set-PriceText "D23" "70.26"
$ws.Range("E23").Value = "  +0.37%  "

# Row 24 - BitcoinCash
Set-PriceText "D24" "264.93"
$ws.Range("E24").Value = "  -0.83%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +1.07%  "

# Row 26 - Filecoin
Set-PriceText "D26" "7.87"
$ws.Range("E26").Value = "  -6.52%  "

# Row 27 - EthereumClassic
Set-PriceText "D27" "26.99"
$ws.Range("E27").Value = "  +2.19%  "

# Row 28 - RenderToken
$ws.Range("E28").Value = "  -2.03%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  -3.97%  "

# Row 31 - now Cosmos (was Hedera)
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-PriceText "D31" "10.72"
$ws.Range("E31").Value = "  +4.30%  "

# Row 32 - now Hedera (was Cosmos)
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-PriceText "D32" "0.105"
$ws.Range("E32").Value = "  -2.74%  "

# Row 33 - VeChain
Set-PriceText "D33" "0.0492"
$ws.Range("E33").Value = "  +10.78%  "

# Row 34 - InjectiveProtocol
Set-PriceText "D34" "36.46"
$ws.Range("E34").Value = "  +6.88%  "

# Row 35 - Toncoin
$ws.Range("E35").Value = "  +0.22%  "

# Row 36 - OKB
Set-PriceText "D36" "49.86"

# Row 37 - FirstDigitalUSD
Set-PriceText "D37" "1.00"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -1.07%  "

# Row 39 - NEARProtocol
Set-PriceText "D39" "4.01"
$ws.Range("E39").Value = "  +9.09%  "

# Row 40 - TheGraph
$ws.Range("E40").Value = "  +0.59%  "

# Row 41 - Monero
Set-PriceText "D41" "129.16"
$ws.Range("E41").Value = "  +0.61%  "

# Row 42 - now Celestia (was ARBITRUM)
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-PriceText "D42" "16.68"
$ws.Range("E42").Value = "  -2.10%  "

# Row 43 - now ARBITRUM (was Celestia)
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-PriceText "D43" "1.84"
$ws.Range("E43").Value = "  -1.22%  "

# Row 44 - Stellar
$ws.Range("E44").Value = "  -0.73%  "

# Row 45 - Stacks
Set-PriceText "D45" "2.52"
$ws.Range("E45").Value = "  -0.29%  "

# Row 46 - EnergySwap
Set-PriceText "D46" "21.75"
$ws.Range("E46").Value = "  +0.06%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  +0.66%  "

# Row 48 - WEMIXToken
$ws.Range("E48").Value = "  -1.82%  "

# Row 49 - Maker
Set-PriceText "D49" "2.071.94"
$ws.Range("E49").Value = "  +1.79%  "

# Row 50 - FlareNetwork
Set-PriceText "D50" "0.0512"
$ws.Range("E50").Value = "  +30.69%  "

# Row 51 - Mantle
Set-PriceText "D51" "0.903"
$ws.Range("E51").Value = "  +10.09%  "
